$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 9478
$ws.Range("C2").Value = 9454
$ws.Range("D2").Value = 8552
$ws.Range("E2").Value = 0.9045906494605458
$ws.Range("F2").Value = 0.9023000633044946
$ws.Range("G2").Value = 0.09509490102641401
$ws.Range("H2").Value = 0.08582945701391566
$ws.Range("I2").Value = 41871748.61682985
$ws.Range("J2").Value = 14736210.06775092
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 14736210.06775092
$ws.Range("M2").Value = 56607958.68458077
$ws.Range("N2").Value = 817193256.2772001
$ws.Range("O2").Value = 797814518.5032001
$ws.Range("P2").Value = 0.01803271130122035
$ws.Range("Q2").Value = 0.0184707218607627

$ws.Range("B3").Value = 9764
$ws.Range("C3").Value = 9740
$ws.Range("D3").Value = 8841
$ws.Range("E3").Value = 0.907700205338809
$ws.Range("F3").Value = 0.9054690700532568
$ws.Range("G3").Value = 0.09318518303774719
$ws.Range("H3").Value = 0.08440088111805848
$ws.Range("I3").Value = 44278388.44033591
$ws.Range("J3").Value = 15669468.00696756
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 15669468.00696756
$ws.Range("M3").Value = 59947856.44730347
$ws.Range("N3").Value = 866217456.3869131
$ws.Range("O3").Value = 846193941.7348431
$ws.Range("P3").Value = 0.01808953155057231
$ws.Range("Q3").Value = 0.01851758472158575

$ws.Range("B4").Value = 10042
$ws.Range("C4").Value = 10027
$ws.Range("D4").Value = 9131
$ws.Range("E4").Value = 0.910641268574848
$ws.Range("F4").Value = 0.9092810197171878
$ws.Range("G4").Value = 0.0918294128663809
$ws.Range("H4").Value = 0.08351367943466678
$ws.Range("I4").Value = 47517624.66176366
$ws.Range("J4").Value = 16969014.67242942
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 16969014.67242942
$ws.Range("M4").Value = 64486639.33419308
$ws.Range("N4").Value = 926735018.4242668
$ws.Range("O4").Value = 906466055.975973
$ws.Range("P4").Value = 0.01831053573575102
$ws.Range("Q4").Value = 0.0187199670197901

$ws.Range("B5").Value = 10325
$ws.Range("C5").Value = 10298
$ws.Range("D5").Value = 9404
$ws.Range("E5").Value = 0.9131870266071082
$ws.Range("F5").Value = 0.9107990314769976
$ws.Range("G5").Value = 0.09067876875608344
$ws.Range("H5").Value = 0.08261628487963281
$ws.Range("I5").Value = 50835217.98122857
$ws.Range("J5").Value = 18294583.82994707
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 18294583.82994707
$ws.Range("M5").Value = 69129801.81117564
$ws.Range("N5").Value = 982952967.4198503
$ws.Range("O5").Value = 962400924.913815
$ws.Range("P5").Value = 0.01861186082785676
$ws.Range("Q5").Value = 0.01900931655025725

$ws.Range("B6").Value = 10639
$ws.Range("C6").Value = 10615
$ws.Range("D6").Value = 9724
$ws.Range("E6").Value = 0.9160621761658031
$ws.Range("F6").Value = 0.9139956762853652
$ws.Range("G6").Value = 0.08959285875955238
$ws.Range("H6").Value = 0.08191004404341455
$ws.Range("I6").Value = 54617108.1983126
$ws.Range("J6").Value = 19814040.63745773
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 19814040.63745773
$ws.Range("M6").Value = 74431148.83577034
$ws.Range("N6").Value = 1044291442.661237
$ws.Range("O6").Value = 1022775193.783701
$ws.Range("P6").Value = 0.01897366944515441
$ws.Range("Q6").Value = 0.01937282088760557
